# Automatische test-sync: 2025-06-29 15:02:50
#
# Adds the 10th test-mail row to the "Logs" sheet, adds the matching
# "IT / Technisch probleem" tally row to the "Dashboard" sheet, extends the
# conditional formatting ranges on "Logs" to cover the new row, and extends
# the bar chart's category/value series references on "Dashboard" to include
# the new data point.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Logs sheet: append row 25
# ---------------------------------------------------------------------------
$logs = $wb.Worksheets.Item("Logs")

$logs.Cells.Item(25, 1).Value = "Ik wil dat je dit regelt met support."
$logs.Cells.Item(25, 2).Value = "mailmind.test@zohomail.eu"
$logs.Cells.Item(25, 3).Value = "Testmail #10: Ik wil dat je dit regelt met support."
$logs.Cells.Item(25, 4).Value = "IT / Technisch probleem"
$logs.Cells.Item(25, 5).Value = "Beste klant,`nBedankt voor uw e-mail. Om uw verzoek efficiënt te kunnen afhandelen, zou ik graag wat meer informatie van u willen ontvangen. Kunt u mij alstublieft meer details geven over het specifieke probleem waarmee u hulp nodig heeft? Op die manier kunnen wij u beter van dienst zijn.`nMet vriendelijke groet,`n[Naam]  `nE-mailassistent  `n[Bedrijfsnaam]"
$logs.Cells.Item(25, 6).Value = "2025-06-29 15:02:17"
$logs.Cells.Item(25, 7).Value = "Ja"
$logs.Cells.Item(25, 8).Value = "Nee"
$logs.Cells.Item(25, 9).Value = "Ja"

# Extend the conditional formatting ranges (D/G/H/I) from row 24 to row 25,
# keeping the same rules/priorities/dxfIds intact.
$logs.Range("D2:D24").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("D2:D25"))
$logs.Range("G2:G24").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("G2:G25"))
$logs.Range("H2:H24").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("H2:H25"))
$logs.Range("I2:I24").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("I2:I25"))

# ---------------------------------------------------------------------------
# 2. Dashboard sheet: append row 8 (tally for the new category)
# ---------------------------------------------------------------------------
$dash = $wb.Worksheets.Item("Dashboard")

$dash.Cells.Item(8, 1).Value = "IT / Technisch probleem"
$dash.Cells.Item(8, 2).Value = 1

# ---------------------------------------------------------------------------
# 3. Dashboard chart: extend the category/value series ranges to row 8
# ---------------------------------------------------------------------------
$chart = $dash.ChartObjects(1).Chart
$series = $chart.SeriesCollection(1)
$series.XValues = "='Dashboard'!`$A`$2:`$A`$8"
$series.Values = "='Dashboard'!`$B`$2:`$B`$8"
